$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 465, shifting the existing rows 465:563 down to 466:564.
$ws.Rows.Item(465).Insert()

# Populate the newly inserted row 465 with a fresh weekly price record for
# "Feria Lagunitas de Puerto Montt" / Cebollín, matching the surrounding rows.
$ws.Range("A465").Value = 4
$ws.Range("B465").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C465").Value = "Los Lagos"
$ws.Range("D465").Value = 45275
$ws.Range("E465").Value = 10
$ws.Range("F465").Value = 100112037
$ws.Range("G465").Value = "Cebollín"
$ws.Range("H465").Value = "Sin especificar"
$ws.Range("I465").Value = "Primera"
$ws.Range("J465").Value = 180
$ws.Range("K465").Value = 6500
$ws.Range("L465").Value = 6500
$ws.Range("M465").Value = 6500
$ws.Range("N465").Value = "$/paquete 36 unidades"
$ws.Range("O465").Value = "Región Metropolitana"
$ws.Range("P465").Value = 181
$ws.Range("Q465").Value = 36
$ws.Range("R465").Value = "Hortaliza"
